$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33 and 34 content swap (ARBITRUM <-> ImmutableX), plus updated price/volume values
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "0.7478"
$ws.Range("E33").Value = "  +1.63%  "

$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "1.173"
$ws.Range("E34").Value = "  +5.74%  "

$priceChanges = @{
    2 = @{ D = "26.987.85"; E = "  +2.08%  " }
    3 = @{ D = "1.843.23"; E = "  +1.83%  " }
    4 = @{ D = "1.008"; E = "  +0.17%  " }
    5 = @{ D = "309.89"; E = "  +1.31%  " }
    6 = @{ E = "  +0.21%  " }
    7 = @{ D = "0.4670"; E = "  +3.83%  " }
    8 = @{ D = "0.3621"; E = "  +1.10%  " }
    9 = @{ D = "0.07132"; E = "  +1.14%  " }
    10 = @{ D = "0.9130"; E = "  +2.66%  " }
    11 = @{ D = "19.54"; E = "  +1.10%  " }
    12 = @{ D = "0.07695"; E = "  -1.21%  " }
    13 = @{ D = "1.809.21"; E = "  +1.21%  " }
    14 = @{ D = "5.274"; E = "  +0.11%  " }
    15 = @{ D = "6.406"; E = "  +1.61%  " }
    16 = @{ D = "88.30"; E = "  +4.30%  " }
    17 = @{ E = "  +0.19%  " }
    18 = @{ D = "0.000008587"; E = "  +0.86%  " }
    19 = @{ E = "  +0.18%  " }
    20 = @{ D = "27.030.93"; E = "  +2.13%  " }
    21 = @{ D = "14.33"; E = "  +1.08%  " }
    22 = @{ D = "5.018"; E = "  +1.15%  " }
    23 = @{ D = "10.64"; E = "  +1.31%  " }
    24 = @{ D = "1.930"; E = "  -1.53%  " }
    25 = @{ D = "152.62"; E = "  +0.66%  " }
    26 = @{ D = "18.23"; E = "  +2.61%  " }
    27 = @{ D = "2.037"; E = "  -0.66%  " }
    28 = @{ D = "114.07"; E = "  +1.76%  " }
    29 = @{ D = "4.898"; E = "  +0.94%  " }
    30 = @{ D = "0.08857"; E = "  +1.99%  " }
    31 = @{ D = "3.197"; E = "  +2.66%  " }
    32 = @{ D = "2.842"; E = "  +3.59%  " }
    35 = @{ D = "4.462"; E = "  +0.49%  " }
    36 = @{ D = "1.081"; E = "  +0.88%  " }
    37 = @{ D = "2.987"; E = "  +3.28%  " }
    38 = @{ D = "0.01940"; E = "  +0.85%  " }
    39 = @{ D = "0.05167"; E = "  +0.86%  " }
    40 = @{ D = "0.5180"; E = "  +1.80%  " }
    41 = @{ D = "6.910"; E = "  +2.17%  " }
    43 = @{ D = "8.127"; E = "  +1.13%  " }
    44 = @{ E = "  +4.99%  " }
    45 = @{ D = "0.4696"; E = "  +0.65%  " }
    46 = @{ E = "  +0.28%  " }
    47 = @{ E = "  +1.09%  " }
    48 = @{ E = "  +2.46%  " }
    49 = @{ D = "0.06045"; E = "  +0.81%  " }
    50 = @{ D = "64.72"; E = "  +1.91%  " }
    51 = @{ D = "36.24"; E = "  +0.97%  " }
}

foreach ($row in $priceChanges.Keys) {
    $entry = $priceChanges[$row]
    if ($entry.ContainsKey("D")) {
        $ws.Range("D$row").Value = $entry.D
    }
    $ws.Range("E$row").Value = $entry.E
}
